$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Oberliga Hamburg")

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-Rows 64 66
Swap-Rows 65 67
Swap-Rows 163 164
Swap-Rows 214 215
